$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "S: Igen, körülményeim nyugodtak, az internet kapcsolat stabil, a feladatokra tudok szánni 60 percet.          D: Nem, nincs lehetőségem nyugodt körülmények között elvégezni a feladatot."
$ws.Range("B6").Value = "D: 2mp    F: 3mp    J: 4mp   K: 5mp"
$ws.Range("D6").Value = "A helyes válasz 4 mp. "

$ws.Range("B5").Select()
